# Nieuwe data toegevoegd via Streamlit op 2024-12-03 18:09:37
# Append a new record row (row 73) to the CompaNanny database sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 73

$ws.Range("A$newRow").Value = "CompaNanny"
$ws.Range("B$newRow").Value = "CompaNanny Prinsenhof"
$ws.Range("C$newRow").Value = "KDV"

# The report date column stores plain text like "YYYY-MM-DD" in this sheet
# (not a real Excel date serial), so force text formatting before writing it
# to avoid COM's automatic date detection, then drop back to the default
# "Normal" style so no extra number-format styling sticks to the cell.
$ws.Range("D$newRow").NumberFormat = "@"
$ws.Range("D$newRow").Value = "2023-12-14"
$ws.Range("D$newRow").Style = "Normal"

$ws.Range("E$newRow").Value = 0
$ws.Range("F$newRow").Value = 0
$ws.Range("G$newRow").Value = 0
$ws.Range("H$newRow").Value = 0
$ws.Range("I$newRow").Value = 0
$ws.Range("J$newRow").Value = 0
